$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1404
$ws1.Range("F4").Value = 25591
$ws1.Range("F5").Value = 564
$ws1.Range("F6").Value = 237
$ws1.Range("F8").Value = 159
$ws1.Range("F9").Value = 413
$ws1.Range("F11").Value = 340
$ws1.Range("F18").Value = 1460
$ws1.Range("F19").Value = 146
$ws1.Range("F20").Value = 404
$ws1.Range("F22").Value = 112

# 本地生活 (Local life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 4864
$ws3.Range("F3").Value = 169

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1404
$ws4.Range("F4").Value = 4864
$ws4.Range("F5").Value = 169
$ws4.Range("F6").Value = 25592
$ws4.Range("F7").Value = 564
$ws4.Range("F9").Value = 237
$ws4.Range("F14").Value = 159
$ws4.Range("F21").Value = 413
$ws4.Range("F24").Value = 340
$ws4.Range("F35").Value = 1460
$ws4.Range("F36").Value = 146
$ws4.Range("F38").Value = 404
$ws4.Range("F40").Value = 112
